$d = $word.ActiveDocument

# Use Find to locate text, then assign Range.Text directly rather than using
# Find.Execute's Replace argument. Find.Execute's built-in replace pushes the
# replacement text through Word's "typing" pipeline, which silently mangles
# plain ASCII apostrophes into curly/typographic quotes (AutoCorrect smart
# quotes). Setting Range.Text after a successful Find leaves the text exactly
# as provided.
function Replace-Text($old, $new) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        throw "Could not find text: $old"
    }
    $r.Text = $new
}

# Title / heading text (appears twice: H1 heading and bold byline near the end)
Replace-Text "Play Fat Santa Slot for Free - Enjoy Festive Gameplay" "Play Fat Santa Free - Festive Christmas Slot Game"
Replace-Text "Play Fat Santa Slot for Free - Enjoy Festive Gameplay" "Play Fat Santa Free - Festive Christmas Slot Game"

# "What we like" bullet list
Replace-Text "Festive atmosphere with cheerful graphics and catchy soundtrack" "Festive Christmas theme"
Replace-Text "Bonus features for more chances to win big" "Randomly triggered Santa's Sleigh feature"
Replace-Text "Compatible with mobile devices" "Free Spins feature with increasing wild symbols"
Replace-Text "Control panel with options to set up to 100 automatic spins" "Cheerful graphics and catchy soundtrack"

# "What we don't like" bullet list
Replace-Text "Lowest paying symbols still pay relatively low" "Limited number of initial free spins"

# Meta description italic text
Replace-Text "Read our review of Fat Santa slots and play for free. Enjoy the festive theme, bonus features, and mobile compatibility of this cheerful game." "Read our review of Fat Santa and play this festive Christmas slot game for free."
